$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.144.41"
$ws.Range("E2").Value = "  +0.75%  "

$ws.Range("D3").Value = "1.875.65"
$ws.Range("E3").Value = "  +0.99%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.32"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.60%  "

$ws.Range("E6").Value = "  +0.06%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5122"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.54%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3919"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +3.24%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08324"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.40%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.119"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.26%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.50"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.12%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.208"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.46%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.65"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.29%  "

$ws.Range("D14").Value = "1.863.30"
$ws.Range("E14").Value = "  +0.31%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.261"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.27%  "

$ws.Range("E16").Value = "  -0.18%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001100"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.83%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.05"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.02%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06639"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.51%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.75"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.46%  "

$ws.Range("E21").Value = "  -0.01%  "

$ws.Range("E22").Value = "  +0.42%  "

$ws.Range("D23").Value = "28.168.93"
$ws.Range("E23").Value = "  +0.72%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.13"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.32%  "

$ws.Range("E25").Value = "  +1.72%  "

$ws.Range("E26").Value = "  -0.32%  "

$ws.Range("D27").Value = "2.090.28"
$ws.Range("E27").Value = "  +1.07%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.500"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.51%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "158.47"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.91%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.59"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.02%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "124.92"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.47%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1064"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.71%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.039"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.26%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.862"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +5.07%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.593"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.27%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.728"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.01%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02454"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.43%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06530"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.48%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2184"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.69%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.203"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.05%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6500"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.81%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.232"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.70%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.987"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.79%  "

$ws.Range("E44").Value = "  +0.58%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6126"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.47%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.01"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.05%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.282"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.08%  "

$ws.Range("E48").Value = "  +0.22%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.007"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.00%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.231"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.22%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "121.22"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.79%  "
